# Add a new "Number of employees / Assets / Turnover" breakdown table
# (by enterprise size: Micro, Small, Medium, Large) below the existing
# MSME summary tables on the Summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: bold header row, same look as the existing "Micro/SMEs/MSMEs"
# header rows (B11:D11 and B17:D17).
$ws.Cells.Item(21, 2).Value = "Number of employees"
$ws.Cells.Item(21, 3).Value = "Assets (local currency, unless noted otherwise)"
$ws.Cells.Item(21, 4).Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B21:D21").Font.Bold = $true

# Rows 22-25: one row per enterprise-size category; data columns are
# left blank (no figures available yet).
$ws.Cells.Item(22, 1).Value = "Micro"
$ws.Cells.Item(22, 2).Value = ""
$ws.Cells.Item(22, 3).Value = ""
$ws.Cells.Item(22, 4).Value = ""

$ws.Cells.Item(23, 1).Value = "Small"
$ws.Cells.Item(23, 2).Value = ""
$ws.Cells.Item(23, 3).Value = ""
$ws.Cells.Item(23, 4).Value = ""

$ws.Cells.Item(24, 1).Value = "Medium"
$ws.Cells.Item(24, 2).Value = ""
$ws.Cells.Item(24, 3).Value = ""
$ws.Cells.Item(24, 4).Value = ""

$ws.Cells.Item(25, 1).Value = "Large"
$ws.Cells.Item(25, 2).Value = ""
$ws.Cells.Item(25, 3).Value = ""
$ws.Cells.Item(25, 4).Value = ""
